$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Jalen Johnson"
$ws.Range("C6").Value = "Atlanta Hawks"

$ws.Range("A8").Value = "Ochai Agbaji"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Toronto Raptors"

$ws.Range("A9").Value = "Karl-Anthony Towns"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "New York Knicks"

$ws.Range("A14").Value = "Keegan Murray"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "Sacramento Kings"

$ws.Range("A15").Value = "Trey Murphy III"
$ws.Range("C15").Value = "New Orleans Pelicans"
